# Applies the ALG.xlsx edit:
#  - Rename groupe labels "2-A"/"2-B"/"2-C" -> "1-A"/"1-B"/"1-C" (shared string rename,
#    affects every student row referencing those groups)
#  - Update each student's "Numero" (date-ish id, column A) from a 2015* value to the
#    matching 2017* value (semester/year re-inscription correction)
#  - Update each student's "Moyenne de l'etudiant" (column E) with the corrected grade

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the 3 group labels everywhere they appear on the sheet ---
$ws.Cells.Replace("2-B", "1-B")
$ws.Cells.Replace("2-C", "1-C")
$ws.Cells.Replace("2-A", "1-A")

# --- Per-row corrections to column A (Numero) and column E (Moyenne) ---
$rows = @(
    @{ Row = 3; A = 20170926; E = 19 },
    @{ Row = 4; A = 20170927; E = 9 },
    @{ Row = 5; A = 20170928; E = 18 },
    @{ Row = 6; A = 20170929; E = 9 },
    @{ Row = 7; A = 20170930; E = 15 },
    @{ Row = 8; A = 20170931; E = 7 },
    @{ Row = 9; A = 20170932; E = 8 },
    @{ Row = 10; A = 20170933; E = 20 },
    @{ Row = 11; A = 20170934; E = 12 },
    @{ Row = 12; A = 20170935; E = 12 },
    @{ Row = 13; A = 20170936; E = 12 },
    @{ Row = 14; A = 20170937; E = 18 },
    @{ Row = 15; A = 20170938; E = 7 },
    @{ Row = 16; A = 20170939; E = 5 },
    @{ Row = 17; A = 20170940; E = 16 },
    @{ Row = 18; A = 20170941; E = 16 },
    @{ Row = 19; A = 20170942; E = 16 },
    @{ Row = 20; A = 20170943; E = 10 },
    @{ Row = 21; A = 20170944; E = 12 },
    @{ Row = 22; A = 20170945; E = 17 },
    @{ Row = 23; A = 20170946; E = 13 },
    @{ Row = 24; A = 20170947; E = 13 },
    @{ Row = 25; A = 20170948; E = 16 },
    @{ Row = 26; A = 20170949; E = 15 },
    @{ Row = 27; A = 20170950; E = 8 },
    @{ Row = 28; A = 20170951; E = 18 },
    @{ Row = 29; A = 20170952; E = 9 },
    @{ Row = 30; A = 20170953; E = 11 },
    @{ Row = 31; A = 20170954; E = 7 },
    @{ Row = 32; A = 20170955; E = 12 },
    @{ Row = 33; A = 20170956 },
    @{ Row = 34; A = 20170957; E = 14 },
    @{ Row = 35; A = 20170958; E = 7 },
    @{ Row = 36; A = 20170959; E = 5 },
    @{ Row = 37; A = 20170960; E = 10 },
    @{ Row = 38; A = 20170961; E = 13 },
    @{ Row = 39; A = 20170962; E = 15 },
    @{ Row = 40; A = 20170963; E = 17 },
    @{ Row = 41; A = 20170964; E = 13 },
    @{ Row = 42; A = 20170965; E = 8 },
    @{ Row = 43; A = 20170966; E = 18 },
    @{ Row = 44; A = 20170967; E = 20 },
    @{ Row = 45; A = 20170968 },
    @{ Row = 46; A = 20170969; E = 12 },
    @{ Row = 47; A = 20170970; E = 7 },
    @{ Row = 48; A = 20170971; E = 10 },
    @{ Row = 49; A = 20170972; E = 9 },
    @{ Row = 50; A = 20170973; E = 12 },
    @{ Row = 51; A = 20170974 },
    @{ Row = 52; A = 20170975; E = 19 },
    @{ Row = 53; A = 20170976; E = 18 },
    @{ Row = 54; A = 20170977; E = 10 },
    @{ Row = 55; A = 20170978; E = 11 },
    @{ Row = 56; A = 20170979 },
    @{ Row = 57; A = 20170980; E = 13 },
    @{ Row = 58; A = 20170981; E = 11 },
    @{ Row = 59; A = 20170982; E = 19 },
    @{ Row = 60; A = 20170983; E = 6 },
    @{ Row = 61; A = 20170984; E = 20 },
    @{ Row = 62; A = 20170985; E = 18 },
    @{ Row = 63; A = 20170986; E = 18 }
)

foreach ($item in $rows) {
    $ws.Cells.Item($item.Row, 1).Value = $item.A
    if ($item.ContainsKey("E")) {
        $ws.Cells.Item($item.Row, 5).Value = $item.E
    }
}
